$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.25
$ws.Range("C2").Value = 0.5
$ws.Range("P2").Value = 0.125
$ws.Range("S2").Value = 0.125
$ws.Range("C3").Value = 0.4285714285714285
$ws.Range("P3").Value = 0.5714285714285714
$ws.Range("P4").Value = 1
$ws.Range("B6").Value = 0.1
$ws.Range("F6").Value = 0.1
$ws.Range("J6").Value = 0.5
$ws.Range("Q6").Value = 0.1
$ws.Range("S6").Value = 0.2
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("J7").Value = 0.3333333333333333
$ws.Range("O7").Value = 0.2222222222222222
$ws.Range("Q7").Value = 0.1111111111111111
$ws.Range("S7").Value = 0.2222222222222222
$ws.Range("D8").Value = 0.04347826086956522
$ws.Range("F8").Value = 0.04347826086956522
$ws.Range("J8").Value = 0.1739130434782609
$ws.Range("Q8").Value = 0.1304347826086956
$ws.Range("R8").Value = 0.08695652173913043
$ws.Range("S8").Value = 0.5217391304347826
$ws.Range("B9").Value = 0.25
$ws.Range("R9").Value = 0.25
$ws.Range("S9").Value = 0.5
$ws.Range("B10").Value = 0.06666666666666667
$ws.Range("D10").Value = 0.04444444444444445
$ws.Range("F10").Value = 0.1111111111111111
$ws.Range("J10").Value = 0.08888888888888889
$ws.Range("Q10").Value = 0.2444444444444444
$ws.Range("R10").Value = 0.06666666666666667
$ws.Range("S10").Value = 0.3777777777777778
$ws.Range("G11").Value = 0.07142857142857142
$ws.Range("J11").Value = 0.07142857142857142
$ws.Range("L11").Value = 0.8571428571428571
$ws.Range("G12").Value = 0.5
$ws.Range("J12").Value = 0.3333333333333333
$ws.Range("S12").Value = 0.1666666666666667
$ws.Range("G13").Value = 1
$ws.Range("H15").Value = 0.1428571428571428
$ws.Range("I15").Value = 0.1428571428571428
$ws.Range("J15").Value = 0.1428571428571428
$ws.Range("S15").Value = 0.5714285714285714
$ws.Range("I16").Value = 0.125
$ws.Range("J16").Value = 0.375
$ws.Range("K16").Value = 0.25
$ws.Range("S16").Value = 0.25
$ws.Range("F17").Value = 0.06666666666666667
$ws.Range("H17").Value = 0.5333333333333333
$ws.Range("J17").Value = 0.2
$ws.Range("K17").Value = 0.06666666666666667
$ws.Range("S17").Value = 0.1333333333333333
$ws.Range("H18").Value = 0.2857142857142857
$ws.Range("I18").Value = 0.1428571428571428
$ws.Range("J18").Value = 0.4285714285714285
$ws.Range("K18").Value = 0.1428571428571428
$ws.Range("F19").Value = 0.01612903225806452
$ws.Range("H19").Value = 0.1935483870967742
$ws.Range("I19").Value = 0.08064516129032258
$ws.Range("J19").Value = 0.2258064516129032
$ws.Range("K19").Value = 0.1612903225806452
$ws.Range("M19").Value = 0.03225806451612903
$ws.Range("O19").Value = 0.06451612903225806
$ws.Range("S19").Value = 0.2258064516129032
